# PROD Suite Update for TC 4,8 and 12
# Replaces the "navigate via mega-menu" flow with a search-box driven flow
# on the TC04_SearchCategory sheet, and updates the corresponding test data
# on the Testdata sheet (category renamed from "Gearing" to
# "Gear Racks and Pinions").

$wb = $excel.ActiveWorkbook
$wsTC = $wb.Worksheets.Item("TC04_SearchCategory")
$wsTD = $wb.Worksheets.Item("Testdata")

# ---------------------------------------------------------------------
# 1) TC04_SearchCategory sheet - insert the extra automation steps
# ---------------------------------------------------------------------

# Make room: 6 new rows just above the old row 5 (old rows 5-7 shift to 11-13)
$wsTC.Range("A5:A10").EntireRow.Insert()
# ... and one more row above what is now row 12 (old row 6/7, now 12/13, shift to 13/14)
$wsTC.Range("A12").EntireRow.Insert()

# Old row 3 (MOUSEOVER / ProductMegamenu / CSS) becomes a simple WAIT step
$wsTC.Range("C3").ClearContents()
$wsTC.Range("D3").ClearContents()
$wsTC.Range("B3").Value = "WAIT"

# Old row 4 (CLICK / Gearing / CSS) becomes a simple WAIT step
$wsTC.Range("C4").ClearContents()
$wsTC.Range("D4").ClearContents()
$wsTC.Range("B4").Value = "WAIT"

# New rows 5-10: search box interaction + waits
$wsTC.Range("B5").Value = "CLICK_PRE_ENTERTEXT"
$wsTC.Range("C5").Value = "SearchBoxHomePage"
$wsTC.Range("D5").Value = "CSS"

$wsTC.Range("B6").Value = "ENTERTEXT"
$wsTC.Range("C6").Value = "SearchBoxHomePage"
$wsTC.Range("D6").Value = "CSS"
$wsTC.Range("E6").Value = "validSearch"

$wsTC.Range("B7").Value = "PRESS_ENTER"
$wsTC.Range("C7").Value = "SearchBoxHomePage"
$wsTC.Range("D7").Value = "CSS"

$wsTC.Range("B8").Value = "WAIT"
$wsTC.Range("B9").Value = "WAIT"
$wsTC.Range("B10").Value = "WAIT"

# Row 11 keeps the original VERIFY_TEXT_PRESENT step (now shifted down)
$wsTC.Range("B11").Value = "VERIFY_TEXT_PRESENT"
$wsTC.Range("C11").Value = "GearingCategoryHeader"
$wsTC.Range("D11").Value = "CSS"
$wsTC.Range("E11").Value = "GearingPLPHeader"

# New row 12: another WAIT before validating the filter
$wsTC.Range("B12").Value = "WAIT"

# Rows 13-14 keep the original trailing verification steps (now shifted down)
$wsTC.Range("B13").Value = "VERIFY_WEBELEMENT_PRESENT"
$wsTC.Range("C13").Value = "CategoryFilter"
$wsTC.Range("D13").Value = "CSS"
$wsTC.Range("E13").Value = "CategoryFilter"

$wsTC.Range("B14").Value = "VERIFY_PAGE_TITLE"
$wsTC.Range("E14").Value = "PLPpageTitle"

# Borders for every data row (1-14) across columns A-E, matching the sheet's
# existing bordered-table look
$wsTC.Range("A2:E14").Borders.LineStyle = 1

$wsTC.Range("B8").Select()

# ---------------------------------------------------------------------
# 2) Testdata sheet - category renamed from "Gearing" to
#    "Gear Racks and Pinions" + new row for the search term
# ---------------------------------------------------------------------

$wsTD.Range("B8").Value = "Gear Racks and Pinions"
$wsTD.Range("B10").Value = "Gear Racks and Pinions | Kaman Industrial"

$wsTD.Range("A17").Value = "validSearch"
$wsTD.Range("B17").Value = "Gearing"
$wsTD.Range("A17:B17").Borders.LineStyle = 1

$wsTD.Columns.Item(2).ColumnWidth = 38

$wsTD.Range("B10").Select()

# Leave TC04_SearchCategory as the active sheet/tab, matching the saved file
$wsTC.Activate()
$wsTC.Range("B8").Select()
